$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-10-11 Friday"; New = "2024-10-12 Saturday" },
    @{ Old = "55×64=3520"; New = "64×23=1472" },
    @{ Old = "51×44=2244"; New = "39×61=2379" },
    @{ Old = "62×43=2666"; New = "87×68=5916" },
    @{ Old = "61×21=1281"; New = "90×36=3240" },
    @{ Old = "99×99=9801"; New = "29×23=667" },
    @{ Old = "33×72=2376"; New = "91×74=6734" },
    @{ Old = "32×95=3040"; New = "45×42=1890" },
    @{ Old = "95×33=3135"; New = "86×32=2752" },
    @{ Old = "92×42=3864"; New = "36×48=1728" },
    @{ Old = "21×66=1386"; New = "60×57=3420" },
    @{ Old = "63×22=1386"; New = "40×98=3920" },
    @{ Old = "99×13=1287"; New = "67×91=6097" },
    @{ Old = "46×11=506"; New = "54×28=1512" },
    @{ Old = "80×28=2240"; New = "81×21=1701" },
    @{ Old = "13×44=572"; New = "19×99=1881" },
    @{ Old = "20×45=900"; New = "85×18=1530" },
    @{ Old = "16×28=448"; New = "97×61=5917" },
    @{ Old = "24×41=984"; New = "52×86=4472" },
    @{ Old = "75×12=900"; New = "91×86=7826" },
    @{ Old = "46×58=2668"; New = "57×78=4446" },
    @{ Old = "24×25=600"; New = "79×99=7821" },
    @{ Old = "49×45=2205"; New = "76×59=4484" },
    @{ Old = "62×38=2356"; New = "26×50=1300" },
    @{ Old = "87×86=7482"; New = "65×48=3120" },
    @{ Old = "43×25=1075"; New = "80×84=6720" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
